$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: height shrinks from 46.5 to 24 (its keyword/appID text is unchanged)
$ws.Rows.Item(12).RowHeight = 24

# Rows 13-16: the keyword/appID list got re-ordered/edited
$ws.Range("A13").Value = "passive income"
$ws.Range("B13").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A14").Value = "travelpayouts"
$ws.Range("B14").Value = "affiliate.marketing.guide"

$ws.Range("A15").Value = "blockchain"
$ws.Range("B15").Value = "block.chain.technology"

$ws.Range("A16").Value = "passive income"
$ws.Range("B16").Value = "affiliate.marketing.guide"

# The sheet shrank by one row: drop the now-superfluous trailing blank row 20
$ws.Rows.Item(20).Delete()

# Match the updated view state: scrolled so row 4 is the top row, A13 selected
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A13").Select()
